# Daily attendance processing - 2025-10-30 11:42:45
#
# Normalizes the "Recorded By" column (G) on the active sheet: whenever the
# comma-separated list of recorders contains the literal token "System"
# but it isn't already first, move it to the front of the list while
# preserving the relative order of the remaining entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

$col = 7  # column G = "Recorded By"

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -isnot [string]) { continue }
    if ($val -eq "") { continue }

    $parts = $val -split ", "
    if ($parts.Count -le 1) { continue }

    $idx = -1
    for ($i = 0; $i -lt $parts.Count; $i++) {
        if ($parts[$i].Equals("System")) {
            $idx = $i
            break
        }
    }

    if ($idx -le 0) { continue }  # no "System" token, or already first

    $newParts = @($parts[$idx])
    for ($i = 0; $i -lt $parts.Count; $i++) {
        if ($i -ne $idx) {
            $newParts += $parts[$i]
        }
    }

    $newVal = [string]::Join(", ", $newParts)
    if ($newVal -ne $val) {
        $cell.Value2 = $newVal
    }
}
